# Update rows with new TVN sampling data (values shifted due to new sampling rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = 91
$ws.Cells.Item(4, 2).Value = "6:11 AM"
$ws.Cells.Item(4, 3).Value = 2725.88
$ws.Cells.Item(4, 4).Value = 2755.88
$ws.Cells.Item(4, 5).Value = 22.16
$ws.Cells.Item(5, 1).Value = 100
$ws.Cells.Item(5, 2).Value = "6:15 AM"
$ws.Cells.Item(5, 3).Value = 2986.3175
$ws.Cells.Item(5, 4).Value = 3016.3175
$ws.Cells.Item(5, 5).Value = 30.9975
$ws.Cells.Item(6, 1).Value = 123
$ws.Cells.Item(6, 2).Value = "6:27 AM"
$ws.Cells.Item(6, 3).Value = 3687.43
$ws.Cells.Item(6, 4).Value = 3717.43
$ws.Cells.Item(6, 5).Value = 25.2
$ws.Cells.Item(7, 1).Value = 151
$ws.Cells.Item(7, 2).Value = "6:40 AM"
$ws.Cells.Item(7, 3).Value = 4506.15
$ws.Cells.Item(7, 4).Value = 4536.15
$ws.Cells.Item(7, 5).Value = 33.81
$ws.Cells.Item(8, 1).Value = 182
$ws.Cells.Item(8, 2).Value = "6:56 AM"
$ws.Cells.Item(8, 3).Value = 5440.89
$ws.Cells.Item(8, 4).Value = 5470.89
$ws.Cells.Item(8, 5).Value = 108.74
$ws.Cells.Item(9, 1).Value = 230
$ws.Cells.Item(9, 2).Value = "7:20 AM"
$ws.Cells.Item(9, 3).Value = 6880.786667
$ws.Cells.Item(9, 4).Value = 6910.786667
$ws.Cells.Item(9, 5).Value = 32.193333
$ws.Cells.Item(10, 1).Value = 239
$ws.Cells.Item(10, 2).Value = "7:24 AM"
$ws.Cells.Item(10, 3).Value = 7146.62
$ws.Cells.Item(10, 4).Value = 7176.62
$ws.Cells.Item(10, 5).Value = 32.95
$ws.Cells.Item(11, 1).Value = 257
$ws.Cells.Item(11, 2).Value = "7:33 AM"
$ws.Cells.Item(11, 3).Value = 7691.65
$ws.Cells.Item(11, 4).Value = 7721.65
$ws.Cells.Item(11, 5).Value = 35.775
$ws.Cells.Item(12, 1).Value = 331
$ws.Cells.Item(12, 2).Value = "8:10 AM"
$ws.Cells.Item(12, 3).Value = 9903.9
$ws.Cells.Item(12, 4).Value = 9933.9
$ws.Cells.Item(12, 5).Value = 66.37
$ws.Cells.Item(13, 1).Value = 524
$ws.Cells.Item(13, 2).Value = "9:47 AM"
$ws.Cells.Item(13, 3).Value = 15702.116667
$ws.Cells.Item(13, 4).Value = 15732.116667
$ws.Cells.Item(13, 5).Value = 35.09
$ws.Cells.Item(14, 1).Value = 550
$ws.Cells.Item(14, 2).Value = "10:00 A"
$ws.Cells.Item(14, 3).Value = 16484
$ws.Cells.Item(14, 4).Value = 16514
$ws.Cells.Item(14, 5).Value = 26.08
$ws.Cells.Item(15, 1).Value = 570
$ws.Cells.Item(15, 2).Value = "10:10 A"
$ws.Cells.Item(15, 3).Value = 17087.03
$ws.Cells.Item(15, 4).Value = 17117.03
$ws.Cells.Item(15, 5).Value = 37.05
$ws.Cells.Item(16, 1).Value = 1016
$ws.Cells.Item(16, 2).Value = "1:53 PM"
$ws.Cells.Item(16, 3).Value = 30454.73
$ws.Cells.Item(16, 4).Value = 30484.73
$ws.Cells.Item(16, 5).Value = 24.88
$ws.Cells.Item(17, 1).Value = 1022
$ws.Cells.Item(17, 2).Value = "1:56 PM"
$ws.Cells.Item(17, 3).Value = 30644.296667
$ws.Cells.Item(17, 4).Value = 30674.296667
$ws.Cells.Item(17, 5).Value = 24.493333
$ws.Cells.Item(18, 1).Value = 1033
$ws.Cells.Item(18, 2).Value = "2:01 PM"
$ws.Cells.Item(18, 3).Value = 30974.42
$ws.Cells.Item(18, 4).Value = 31004.42
$ws.Cells.Item(18, 5).Value = 36.19
$ws.Cells.Item(19, 1).Value = 1086
$ws.Cells.Item(19, 2).Value = "2:28 PM"
$ws.Cells.Item(19, 3).Value = 32575.77
$ws.Cells.Item(19, 4).Value = 32605.77
$ws.Cells.Item(19, 5).Value = 21.33
$ws.Cells.Item(20, 1).Value = 1103
$ws.Cells.Item(20, 2).Value = "2:36 PM"
$ws.Cells.Item(20, 3).Value = 33070.88
$ws.Cells.Item(20, 4).Value = 33100.88
$ws.Cells.Item(20, 5).Value = 25.1
$ws.Cells.Item(21, 1).Value = 1141
$ws.Cells.Item(21, 2).Value = "2:55 PM"
$ws.Cells.Item(21, 3).Value = 34219.87
$ws.Cells.Item(21, 4).Value = 34249.87
$ws.Cells.Item(21, 5).Value = 21.48
$ws.Cells.Item(22, 1).Value = 1182
$ws.Cells.Item(22, 2).Value = "3:16 PM"
$ws.Cells.Item(22, 3).Value = 35455.655
$ws.Cells.Item(22, 4).Value = 35485.655
$ws.Cells.Item(22, 5).Value = 37.115
$ws.Cells.Item(23, 1).Value = 1200
$ws.Cells.Item(23, 2).Value = "3:25 PM"
$ws.Cells.Item(23, 3).Value = 35982.38
$ws.Cells.Item(23, 4).Value = 36012.38
$ws.Cells.Item(23, 5).Value = 37.86
$ws.Cells.Item(24, 1).Value = 1214
$ws.Cells.Item(24, 2).Value = "3:32 PM"
$ws.Cells.Item(24, 3).Value = 36391.61
$ws.Cells.Item(24, 4).Value = 36421.61
$ws.Cells.Item(24, 5).Value = 96.52
$ws.Cells.Item(25, 1).Value = 1221
$ws.Cells.Item(25, 2).Value = "3:35 PM"
$ws.Cells.Item(25, 3).Value = 36616.4
$ws.Cells.Item(25, 4).Value = 36646.4
$ws.Cells.Item(25, 5).Value = 23.47
$ws.Cells.Item(26, 1).Value = 1227
$ws.Cells.Item(26, 2).Value = "3:39 PM"
$ws.Cells.Item(26, 3).Value = 36808.87
$ws.Cells.Item(26, 4).Value = 36838.87
$ws.Cells.Item(26, 5).Value = 22.12
$ws.Cells.Item(27, 1).Value = 1245
$ws.Cells.Item(27, 2).Value = "3:47 PM"
$ws.Cells.Item(27, 3).Value = 37332.09
$ws.Cells.Item(27, 4).Value = 37362.09
$ws.Cells.Item(27, 5).Value = 193.12
$ws.Cells.Item(28, 1).Value = 1257
$ws.Cells.Item(28, 2).Value = "3:53 PM"
$ws.Cells.Item(28, 3).Value = 37684.8
$ws.Cells.Item(28, 4).Value = 37714.8
$ws.Cells.Item(28, 5).Value = 36.82
$ws.Cells.Item(30, 1).Value = 1277
$ws.Cells.Item(30, 2).Value = "4:03 PM"
$ws.Cells.Item(30, 3).Value = 38285.38
$ws.Cells.Item(30, 4).Value = 38315.38
$ws.Cells.Item(30, 5).Value = 25.395
$ws.Cells.Item(31, 1).Value = 1283
$ws.Cells.Item(31, 2).Value = "4:06 PM"
$ws.Cells.Item(31, 3).Value = 38479.19
$ws.Cells.Item(31, 4).Value = 38509.19
$ws.Cells.Item(31, 5).Value = 51.24
$ws.Cells.Item(32, 1).Value = 1297
$ws.Cells.Item(32, 2).Value = "4:13 PM"
$ws.Cells.Item(32, 3).Value = 38885.46
$ws.Cells.Item(32, 4).Value = 38915.46
$ws.Cells.Item(32, 5).Value = 73.44
$ws.Cells.Item(33, 1).Value = 1311
$ws.Cells.Item(33, 2).Value = "4:20 PM"
$ws.Cells.Item(33, 3).Value = 39315.5
$ws.Cells.Item(33, 4).Value = 39345.5
$ws.Cells.Item(33, 5).Value = 76.02
$ws.Cells.Item(34, 1).Value = 1323
$ws.Cells.Item(34, 2).Value = "4:26 PM"
$ws.Cells.Item(34, 3).Value = 39677.255
$ws.Cells.Item(34, 4).Value = 39707.255
$ws.Cells.Item(34, 5).Value = 57.735
$ws.Cells.Item(35, 1).Value = 1335
$ws.Cells.Item(35, 2).Value = "4:32 PM"
$ws.Cells.Item(35, 3).Value = 40027.72
$ws.Cells.Item(35, 4).Value = 40057.72
$ws.Cells.Item(35, 5).Value = 47.705
$ws.Cells.Item(36, 1).Value = 1341
$ws.Cells.Item(36, 2).Value = "4:35 PM"
$ws.Cells.Item(36, 3).Value = 40209.895
$ws.Cells.Item(36, 4).Value = 40239.895
$ws.Cells.Item(36, 5).Value = 54.99
$ws.Cells.Item(37, 1).Value = 1349
$ws.Cells.Item(37, 2).Value = "4:39 PM"
$ws.Cells.Item(37, 3).Value = 40446.58
$ws.Cells.Item(37, 4).Value = 40476.58
$ws.Cells.Item(37, 5).Value = 52.82
$ws.Cells.Item(38, 1).Value = 1360
$ws.Cells.Item(38, 2).Value = "4:45 PM"
$ws.Cells.Item(38, 3).Value = 40775.5
$ws.Cells.Item(38, 4).Value = 40805.5
$ws.Cells.Item(38, 5).Value = 27.96
$ws.Cells.Item(39, 1).Value = 1372
$ws.Cells.Item(39, 2).Value = "4:51 PM"
$ws.Cells.Item(39, 3).Value = 41141.465
$ws.Cells.Item(39, 4).Value = 41171.465
$ws.Cells.Item(39, 5).Value = 98.08
$ws.Cells.Item(41, 1).Value = 1412
$ws.Cells.Item(41, 2).Value = "5:11 PM"
$ws.Cells.Item(41, 3).Value = 42334.07
$ws.Cells.Item(41, 4).Value = 42364.07
$ws.Cells.Item(41, 5).Value = 101.66
$ws.Cells.Item(42, 1).Value = 1418
$ws.Cells.Item(42, 2).Value = "5:14 PM"
$ws.Cells.Item(42, 3).Value = 42521.595
$ws.Cells.Item(42, 4).Value = 42551.595
$ws.Cells.Item(42, 5).Value = 25.965
$ws.Cells.Item(43, 1).Value = 1432
$ws.Cells.Item(43, 2).Value = "5:21 PM"
$ws.Cells.Item(43, 3).Value = 42943.82
$ws.Cells.Item(43, 4).Value = 42973.82
$ws.Cells.Item(43, 5).Value = 31.495
